$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 2297
$ws.Range("I9").Value = 2970.5
$ws.Range("K9").Value = 2970.5
$ws.Range("M9").Value = -2801.5
# Row 11
$ws.Range("H11").Value = 64909.06
$ws.Range("I11").Value = 64909.06
$ws.Range("K11").Value = 64909.06
$ws.Range("M11").Value = -64769.06
# Row 12
$ws.Range("H12").Value = 203.14285
$ws.Range("I12").Value = 199.5
$ws.Range("K12").Value = 199.5
$ws.Range("M12").Value = -29.5
# Row 43
$ws.Range("H43").Value = 1075
$ws.Range("J43").Value = 1137.5
$ws.Range("L43").Value = 1137.5
$ws.Range("N43").Value = -1275.5
# Row 64
$ws.Range("H64").Value = 4249.75
$ws.Range("J64").Value = 4285.5713
$ws.Range("L64").Value = 4285.5713
$ws.Range("N64").Value = -4781.5713
# Row 67
$ws.Range("H67").Value = 4249.75
$ws.Range("J67").Value = 4285.5713
$ws.Range("L67").Value = 4285.5713
$ws.Range("N67").Value = -6001.5713
# Row 132
$ws.Range("H132").Value = 3094.9363
$ws.Range("I132").Value = 2819.932
$ws.Range("J132").Value = 7128.3335
$ws.Range("K132").Value = 8459.795999999998
$ws.Range("L132").Value = 21385.0005
$ws.Range("M132").Value = -5929.795999999998
$ws.Range("N132").Value = -26445.0005

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 13336413
$ws.Range("I61").Value = 16669415
$ws.Range("K61").Value = 16669415
$ws.Range("M61").Value = -16669203
# Row 74
$ws.Range("H74").Value = 2558.5908
$ws.Range("I74").Value = 2537.4375
$ws.Range("J74").Value = 2615
$ws.Range("K74").Value = 2537.4375
$ws.Range("L74").Value = 2615
$ws.Range("M74").Value = -1663.4375
$ws.Range("N74").Value = -4363
# Row 77
$ws.Range("H77").Value = 2558.5908
$ws.Range("I77").Value = 2537.4375
$ws.Range("J77").Value = 2615
$ws.Range("K77").Value = 12687.1875
$ws.Range("L77").Value = 13075
$ws.Range("M77").Value = -8319.1875
$ws.Range("N77").Value = -21811
# Row 132
$ws.Range("H132").Value = 21278056
$ws.Range("I132").Value = 22223630
$ws.Range("K132").Value = 66670890
$ws.Range("M132").Value = -66668360
# Row 136
$ws.Range("H136").Value = 13336413
$ws.Range("I136").Value = 16669415
$ws.Range("K136").Value = 50008245
$ws.Range("M136").Value = -50005695

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 1583.7693
$ws.Range("I105").Value = 1488.9546
$ws.Range("K105").Value = 1488.9546
$ws.Range("M105").Value = 258.0454
# Row 134
$ws.Range("H134").Value = 1726.4
$ws.Range("J134").Value = 1798.25
$ws.Range("L134").Value = 5394.75
$ws.Range("N134").Value = -10464.75

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 417.2
$ws.Range("I22").Value = 295.5
$ws.Range("J22").Value = 498.33334
$ws.Range("K22").Value = 295.5
$ws.Range("L22").Value = 498.33334
$ws.Range("M22").Value = 54.5
$ws.Range("N22").Value = -1198.33334
# Row 31
$ws.Range("H31").Value = 2832.7432
$ws.Range("I31").Value = 1740.1428
$ws.Range("K31").Value = 1740.1428
$ws.Range("M31").Value = -1445.1428
# Row 34
$ws.Range("H34").Value = 2832.7432
$ws.Range("I34").Value = 1740.1428
$ws.Range("K34").Value = 1740.1428
$ws.Range("M34").Value = -1538.1428
# Row 62
$ws.Range("H62").Value = 76930000
$ws.Range("J62").Value = 250006450
$ws.Range("L62").Value = 250006450
$ws.Range("N62").Value = -250007698
# Row 65
$ws.Range("H65").Value = 76930000
$ws.Range("J65").Value = 250006450
$ws.Range("L65").Value = 1250032250
$ws.Range("N65").Value = -1250038490
# Row 103
$ws.Range("H103").Value = 9684.223
$ws.Range("I103").Value = 9684.223
$ws.Range("K103").Value = 9684.223
$ws.Range("M103").Value = -8512.223
# Row 134
$ws.Range("H134").Value = 2198.25
$ws.Range("I134").Value = 1235.25
$ws.Range("K134").Value = 3705.75
$ws.Range("M134").Value = -1170.75

$ws = $wb.Worksheets.Item("CUL")
# Row 97
$ws.Range("H97").Value = 1349.4286
$ws.Range("I97").Value = 1688.25
$ws.Range("J97").Value = 897.6667
$ws.Range("K97").Value = 5064.75
$ws.Range("L97").Value = 2693.0001
$ws.Range("M97").Value = -4568.75
$ws.Range("N97").Value = -3685.0001
# Row 136
$ws.Range("H136").Value = 2805.5557
$ws.Range("I136").Value = 9250
$ws.Range("K136").Value = 27750
$ws.Range("M136").Value = -22650
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()
# Row 139
$ws.Range("H139").Value = 8335683.5
$ws.Range("I139").Value = 16667617
$ws.Range("K139").Value = 50002851
$ws.Range("M139").Value = -49997711
# Row 140
$ws.Range("H140").Value = 1387.0952
$ws.Range("J140").Value = 2940
$ws.Range("L140").Value = 8820
$ws.Range("N140").Value = -19180

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4493.75
$ws.Range("I40").Value = 3316.6667
$ws.Range("K40").Value = 3316.6667
$ws.Range("M40").Value = -3180.6667
# Row 43
$ws.Range("H43").Value = 6265833.5
$ws.Range("I43").Value = 3598571.5
$ws.Range("K43").Value = 3598571.5
$ws.Range("M43").Value = -3598378.5
# Row 46
$ws.Range("H46").Value = 1675.4
$ws.Range("I46").Value = 907.6667
$ws.Range("J46").Value = 1810.8823
$ws.Range("K46").Value = 907.6667
$ws.Range("L46").Value = 1810.8823
$ws.Range("M46").Value = -719.6667
$ws.Range("N46").Value = -2186.8823
# Row 55
$ws.Range("H55").Value = 641.6667
$ws.Range("I55").Value = 565.5
$ws.Range("J55").Value = 794
$ws.Range("K55").Value = 565.5
$ws.Range("L55").Value = 794
$ws.Range("M55").Value = -392.5
$ws.Range("N55").Value = -1140
# Row 132
$ws.Range("H132").Value = 5999.32
$ws.Range("J132").Value = 7592.2
$ws.Range("L132").Value = 22776.6
$ws.Range("N132").Value = -27836.6
# Row 136
$ws.Range("H136").Value = 2241.7368
$ws.Range("I136").Value = 2142.6428
$ws.Range("J136").Value = 2519.2
$ws.Range("K136").Value = 6427.928400000001
$ws.Range("L136").Value = 7557.599999999999
$ws.Range("M136").Value = -3877.928400000001
$ws.Range("N136").Value = -12657.6

$ws = $wb.Worksheets.Item("WVR")
# Row 64
$ws.Range("H64").Value = 53777
$ws.Range("I64").Value = 53777
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 53777
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -53529
$ws.Range("N64").ClearContents()
# Row 67
$ws.Range("H67").Value = 53777
$ws.Range("I67").Value = 53777
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 53777
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -52919
$ws.Range("N67").ClearContents()
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
# Row 132
$ws.Range("H132").Value = 4275
$ws.Range("I132").Value = 4110.393
$ws.Range("J132").Value = 4787.1113
$ws.Range("K132").Value = 12331.179
$ws.Range("L132").Value = 14361.3339
$ws.Range("M132").Value = -9801.179
$ws.Range("N132").Value = -19421.3339
# Row 136
$ws.Range("H136").Value = 2853.4075
$ws.Range("I136").Value = 1097.619
$ws.Range("K136").Value = 3292.857
$ws.Range("M136").Value = -742.857
